$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Non-numeric text columns (Coin name, Link, Volume label) ---
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E7').Value = '6MXTokenMX'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E8').Value = '7FTXTokenFTT'
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('E9').Value = '8WazirXWRX'
$ws.Range('B10').Value = 'MandalaExchangeToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('E10').Value = '9MandalaExchangeTokenMDX'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('E11').Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('E12').Value = '11BitrueCoinBTR'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('E13').Value = '12BitMartTokenBMX'
$ws.Range('B14').Value = 'MCDex'
$ws.Range('C14').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('E14').Value = '13MCDexMCB'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('E15').Value = '14BitForexTokenBF'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('E16').Value = '15CoinExTokenCET'
$ws.Range('B17').Value = 'One'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E17').Value = '16OneONE'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('E18').Value = '17TigerCashTCH'
$ws.Range('B19').Value = 'HotbitToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('E19').Value = '18HotbitTokenHTB'
$ws.Range('B20').Value = 'BitKan'
$ws.Range('C20').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('E20').Value = '19BitKanKAN'
$ws.Range('B21').Value = 'NitroEx'
$ws.Range('C21').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('E21').Value = '20NitroExNTX'
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E22').Value = '21LEOLEO'
$ws.Range('B23').Value = 'KuCoinToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('E23').Value = '22KuCoinTokenKCS'
$ws.Range('E43').Value = '42CEJICEJIBestin24h'

# --- Price column (D): numeric-looking strings must stay as text ---
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.27'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.448'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05742'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.437'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8117'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8925'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1443'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07379'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.02991'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03088'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09409'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.928'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001586'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.04827'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0005851'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.006217'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.005109'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0009954'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0001500'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.747'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.328'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.3278'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03914'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006788'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1070'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.003200'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007384'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005640'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1690'
$ws.Range('D48').Style = 'Normal'
